# Insert a new data row for 2026/02/01 (日, 16:00 slot) into the Sheet1
# time-series table. This shifts the existing row 752 ("2026/12/29" block)
# and everything below it down by one row (752 -> 753, ..., 793 -> 794),
# matching the new date/ranking entry that was appended upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 752 (and everything after it) down one row, opening up a blank
# row 752 for the new entry.
$ws.Rows("752").Insert()

# Fill the newly-opened row. The leading apostrophe forces column A to be
# stored as literal text (matching the rest of the sheet) instead of being
# auto-converted to a date serial number by the General number format.
# Column B ("日") is plain text already and needs no such guard.
$ws.Range("A752").Value = "'2026/02/01"
$ws.Range("B752").Value = "日"
$ws.Range("C752").Value = 16
$ws.Range("D752").Value = 201
